# Applies the 2023 Sweden Allsvenskan update:
#  1) A number of match rows had their F:V data (everything except the
#     leading Indice/pais/torneio/temporada/data_partida columns) rotated
#     among themselves (the underlying scrape re-ordered/re-matched rows).
#  2) Two brand-new match rows (226 and 227) were appended at the end.
#  3) The used-range dimension grows accordingly (Excel maintains this
#     automatically as new cells are written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row groups whose F:V block must be rotated "up by one" (the first
#    row receives what used to be in the second row, and so on, with the
#    last row in the group wrapping around to receive what used to be in
#    the first row). Two-row groups are therefore a plain swap.
# ---------------------------------------------------------------------
$groups = @(
    ,@(22,23)
    ,@(34,35)
    ,@(38,39)
    ,@(48,49)
    ,@(67,68,69)
    ,@(70,71)
    ,@(98,99)
    ,@(114,115)
    ,@(117,118)
    ,@(119,120)
    ,@(127,128)
    ,@(141,142)
    ,@(144,145)
    ,@(149,150,151)
    ,@(157,158)
    ,@(196,197)
    ,@(198,199)
    ,@(205,206,207)
)

foreach ($grp in $groups) {
    $n = $grp.Count

    # Snapshot the original F:V values for every row in this group first,
    # since we must not read an already-overwritten row.
    $snapshots = @()
    for ($i = 0; $i -lt $n; $i++) {
        $r = $grp[$i]
        $snapshots += , ($ws.Range("F$r`:V$r").Value2)
    }

    # Row i gets what row (i+1) used to hold; the last row wraps to row 0.
    for ($i = 0; $i -lt $n; $i++) {
        $r = $grp[$i]
        $src = $snapshots[($i + 1) % $n]
        $ws.Range("F$r`:V$r").Value = $src
    }
}

# ---------------------------------------------------------------------
# 2) Append the two new rows (226 and 227) at the end of the sheet.
#    Copy row 225 first so the new rows inherit the same per-column
#    styles (s="1" on Indice, s="2" on the date column, etc.), then
#    overwrite the values.
# ---------------------------------------------------------------------
$ws.Range("A225:V225").Copy($ws.Range("A226:V226"))
$ws.Range("A225:V225").Copy($ws.Range("A227:V227"))

# NOTE: a flat PowerShell literal array (@(...)) assigned to a multi-cell
# Range.Value does NOT get distributed one-value-per-cell the way a true
# 2-D SAFEARRAY (what Range.Value2 returns for a multi-cell read) does, so
# each new row is filled in one cell at a time via Cells.Item instead.

$row226 = @(
    225, "sweden", "allsvenskan", "2023", 45234.625,
    "Varnamo", 0, "Hammarby", 0,
    2.89, "30/10/2023 19:13", 2.4, "04/11/2023 14:54",
    3.43, "30/10/2023 19:13", 3.66, "04/11/2023 14:54",
    2.54, "30/10/2023 19:13", 2.96, "04/11/2023 14:54",
    "https://www.betexplorer.com/football/sweden/allsvenskan/varnamo-hammarby/v3VYM2r8/"
)
for ($i = 0; $i -lt $row226.Count; $i++) {
    $ws.Cells.Item(226, 1 + $i).Value = $row226[$i]
}

$row227 = @(
    226, "sweden", "allsvenskan", "2023", 45234.72916666666,
    "Djurgarden", 2, "Sirius", 4,
    1.64, "30/10/2023 19:13", 1.6, "04/11/2023 17:26",
    4.45, "30/10/2023 19:13", 4.71, "04/11/2023 17:26",
    4.83, "30/10/2023 19:13", 5.17, "04/11/2023 17:26",
    "https://www.betexplorer.com/football/sweden/allsvenskan/djurgarden-sirius/2JX1K9M9/"
)
for ($i = 0; $i -lt $row227.Count; $i++) {
    $ws.Cells.Item(227, 1 + $i).Value = $row227[$i]
}

# ---------------------------------------------------------------------
# 3) Make sure the worksheet's recorded dimension covers the new rows.
# ---------------------------------------------------------------------
$ws.Range("A1:V227").Select() | Out-Null

Write-Output "edit applied"
